$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: App -> Course Revenue ---
$ws.Range("A2").Value = "Course Revenue"
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 45726.250231481485

# --- Row 3: Business profit -> App development ---
$ws.Range("A3").Value = "App development"
$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = 45723.250231481485

# --- Row 4 (new): Fizza Selling ---
$ws.Range("A4").Value = "Fizza Selling"
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 45722.250231481485

# --- Row 5 (new): Business income ---
$ws.Range("A5").Value = "Business income"
$ws.Range("B5").Value = 250
$ws.Range("C5").Value = 45717.250231481485

# Copy the date cell formatting (style index) from the existing C2 cell
# down onto the two newly created date cells, so C4/C5 share the same
# numFmt (m/d/yyyy serial date style) as C2/C3 instead of getting the
# generic/default style.
$ws.Range("C2").Copy()
$ws.Range("C4:C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
